$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
questions = [
    {
        "title": "MongoDB has a sharded setup. Application logs show that it cannot connect to the MongoDB database.  Which of the MongoDB components would you check first to ensure that the application can connect to the database?",
        "ques_type": 2,
        "options": [
            "MongoDB config server ",
            "MongoDB mongos",
            "MongoDB primary",
            "MongoDB secondary"
        ],
        "score": "MongoDB mongos"
    },
    {
        "title": "Internal auditing showed that anyone with network access could access MongoDB data in a standalone mongod instance.  What actions would you perform in this situation in order to restrict access to the MongoDB data?",
        "ques_type": 2,
        "options": [
            "Do nothing, as it is only internal access.",
            "Create different roles and users for each database.",
            "Enable access control and create database users with certain roles.",
            "Enable access control and create necessary user roles."
        ],
        "score": "Enable access control and create database users with certain roles."
    },
    {
        "title": "MongoDB logs showed many entries with \u201cCannot create new thread, closing connection\u201d error statements.  What actions would you perform to troubleshoot this issue?",
        "ques_type": 2,
        "options": [
            "Analyze db.currentOp() command output to see slow queries.",
            "Check ulimit settings in the operating system and modify if necessary.",
            "Increase the number of connections to MongoDB from the application side.",
            "Check the operating system CPU and memory utilization and increase resources."
        ],
        "score": "Check ulimit settings in the operating system and modify if necessary."
    },
    {
        "title": "The secondary member of a MongoDB replica set is lagging behind the primary, and the lagging duration keeps increasing.  True or false: In this particular case, it is possible that the secondary cannot connect to the primary even though the primary can connect to the secondary.",
        "ques_type": 11,
        "options": [
            "true",
            "false"
        ],
        "score": "True"
    }
]
'@

# Strip the trailing newline that the here-string literal adds after the closing ]
$newText = $newText.TrimEnd("`r", "`n")

# Row 1 (A1) previously held a bold/bordered placeholder value (0); row 2 (A2) held
# the real shared-string text. Drop the old A1 formatting/value, write the text into
# A1 with the default (unstyled) look, then remove the now-empty second row so the
# sheet is back down to a single A1 cell.
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $newText
$ws.Rows(2).Delete()
